$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 'aa'
$ws.Range("J4").Value = 'Agree/Accept'
$ws.Range("I8").Value = '%'
$ws.Range("J8").Value = 'Uninterpretable'
$ws.Range("I12").Value = 'b'
$ws.Range("J12").Value = 'Acknowledge (Backchannel)'
$ws.Range("I16").Value = 'sd'
$ws.Range("J16").Value = 'Statement-non-opinion'
$ws.Range("I18").Value = 'sd'
$ws.Range("J18").Value = 'Statement-non-opinion'
$ws.Range("I25").Value = 'aa'
$ws.Range("J25").Value = 'Agree/Accept'
$ws.Range("I53").Value = '%'
$ws.Range("J53").Value = 'Uninterpretable'
$ws.Range("I56").Value = 'aa'
$ws.Range("J56").Value = 'Agree/Accept'
$ws.Range("I58").Value = 'b'
$ws.Range("J58").Value = 'Acknowledge (Backchannel)'
$ws.Range("I60").Value = 'aa'
$ws.Range("J60").Value = 'Agree/Accept'
$ws.Range("I61").Value = 'sd'
$ws.Range("J61").Value = 'Statement-non-opinion'
$ws.Range("I63").Value = '%'
$ws.Range("J63").Value = 'Uninterpretable'
$ws.Range("I70").Value = '%'
$ws.Range("J70").Value = 'Uninterpretable'
$ws.Range("I73").Value = '%'
$ws.Range("J73").Value = 'Uninterpretable'
$ws.Range("I92").Value = 'sd'
$ws.Range("J92").Value = 'Statement-non-opinion'
$ws.Range("I97").Value = '%'
$ws.Range("J97").Value = 'Uninterpretable'
$ws.Range("I98").Value = 'aa'
$ws.Range("J98").Value = 'Agree/Accept'
$ws.Range("I104").Value = 'b'
$ws.Range("J104").Value = 'Acknowledge (Backchannel)'
$ws.Range("I106").Value = 'aa'
$ws.Range("J106").Value = 'Agree/Accept'
$ws.Range("I109").Value = 'aa'
$ws.Range("J109").Value = 'Agree/Accept'
$ws.Range("I117").Value = 'aa'
$ws.Range("J117").Value = 'Agree/Accept'
$ws.Range("I129").Value = 'aa'
$ws.Range("J129").Value = 'Agree/Accept'
$ws.Range("I132").Value = 'aa'
$ws.Range("J132").Value = 'Agree/Accept'
$ws.Range("I140").Value = 'sd'
$ws.Range("J140").Value = 'Statement-non-opinion'
$ws.Range("I143").Value = 'sd'
$ws.Range("J143").Value = 'Statement-non-opinion'
$ws.Range("I150").Value = 'b'
$ws.Range("J150").Value = 'Acknowledge (Backchannel)'
$ws.Range("I154").Value = 'ba'
$ws.Range("J154").Value = 'Appreciation'
$ws.Range("I169").Value = 'sd'
$ws.Range("J169").Value = 'Statement-non-opinion'
$ws.Range("I170").Value = '%'
$ws.Range("J170").Value = 'Uninterpretable'
$ws.Range("I187").Value = 'ba'
$ws.Range("J187").Value = 'Appreciation'
$ws.Range("I189").Value = 'sd'
$ws.Range("J189").Value = 'Statement-non-opinion'
$ws.Range("I205").Value = '%'
$ws.Range("J205").Value = 'Uninterpretable'
$ws.Range("I210").Value = 'aa'
$ws.Range("J210").Value = 'Agree/Accept'
$ws.Range("I220").Value = 'b'
$ws.Range("J220").Value = 'Acknowledge (Backchannel)'
$ws.Range("I228").Value = 'aa'
$ws.Range("J228").Value = 'Agree/Accept'
$ws.Range("I229").Value = 'sd'
$ws.Range("J229").Value = 'Statement-non-opinion'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I239").Value = 'aa'
$ws.Range("J239").Value = 'Agree/Accept'
$ws.Range("I242").Value = 'sv'
$ws.Range("J242").Value = 'Statement-opinion'
$ws.Range("I261").Value = 'aa'
$ws.Range("J261").Value = 'Agree/Accept'
$ws.Range("I262").Value = 'aa'
$ws.Range("J262").Value = 'Agree/Accept'
$ws.Range("I270").Value = 'aa'
$ws.Range("J270").Value = 'Agree/Accept'
$ws.Range("I278").Value = '%'
$ws.Range("J278").Value = 'Uninterpretable'
$ws.Range("I322").Value = 'ba'
$ws.Range("J322").Value = 'Appreciation'
$ws.Range("I325").Value = 'sd'
$ws.Range("J325").Value = 'Statement-non-opinion'
$ws.Range("I332").Value = 'sv'
$ws.Range("J332").Value = 'Statement-opinion'
$ws.Range("I333").Value = 'aa'
$ws.Range("J333").Value = 'Agree/Accept'
$ws.Range("I351").Value = 'ba'
$ws.Range("J351").Value = 'Appreciation'
$ws.Range("I358").Value = 'sv'
$ws.Range("J358").Value = 'Statement-opinion'
$ws.Range("I364").Value = 'sd'
$ws.Range("J364").Value = 'Statement-non-opinion'
